# The NATMI export script was re-run with the newer TPM-based pipeline.
# For the Wnt3->Fzd8 pair this collapses the table from 7 data rows (two
# "ligand detection" groups, Sending=ECs and Sending=FAPs) down to 3 rows
# -- only the Sending=FAPs group survives -- and refreshes every
# specificity/weight metric for the rows that remain.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old Sending=ECs block (rows 2-4) worth of rows so that the
# former Sending=FAPs rows (5-7) end up as the new rows 2-4.
$ws.Rows("2:4").Delete() | Out-Null

# Updated data, one row per remaining target cluster (ECs / FAPs / MuSCs).
$data = @(
    @("FAPs", "Wnt3", "Fzd8", "ECs",   2, 0.6666666666666666, 0.1741663333333333, 0.5224989999999999, 1, 1, 3, 1, 3.235341333333333, 9.706023999999999, 0.2153734454473681, 0.2153734454473681, 0.5634875371084443, 5.071387833975999, 0.2153734454473681, 0.2153734454473681),
    @("FAPs", "Wnt3", "Fzd8", "FAPs",  2, 0.6666666666666666, 0.1741663333333333, 0.5224989999999999, 1, 1, 3, 1, 7.273511666666667, 21.820535,          0.4841904166376352, 0.4841904166376352, 1.266800857440555,  11.401207716965,   0.4841904166376352, 0.4841904166376352),
    @("FAPs", "Wnt3", "Fzd8", "MuSCs", 2, 0.6666666666666666, 0.1741663333333333, 0.5224989999999999, 1, 1, 3, 1, 4.513153666666667, 13.539461,          0.3004361379149967, 0.3004361379149967, 0.7860394258932222, 7.074354833038999, 0.3004361379149967, 0.3004361379149967)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($c = 1; $c -le $row.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}

Write-Output "Wnt3-Fzd8: rewrote rows 2-4 with updated TPM specificity values"
